$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.505.32'
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').Value = '1.913.96'
$ws.Range('E3').Value = '  -1.18%  '

$ws.Range('D4').Value = '0.9980'
$ws.Range('E4').Value = '  -0.32%  '

$ws.Range('D5').Value = '239.96'
$ws.Range('E5').Value = '  -1.44%  '

$ws.Range('D6').Value = '0.9983'
$ws.Range('E6').Value = '  -0.24%  '

$ws.Range('D7').Value = '0.4781'
$ws.Range('E7').Value = '  -2.38%  '

$ws.Range('D8').Value = '0.2844'
$ws.Range('E8').Value = '  -3.70%  '

$ws.Range('D9').Value = '0.06717'
$ws.Range('E9').Value = '  -2.51%  '

$ws.Range('D10').Value = '19.44'
$ws.Range('E10').Value = '  +0.73%  '

$ws.Range('D11').Value = '104.26'
$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').Value = '0.07765'
$ws.Range('E12').Value = '  -0.33%  '

$ws.Range('D13').Value = '1.913.85'
$ws.Range('E13').Value = '  -1.16%  '

$ws.Range('D14').Value = '5.228'
$ws.Range('E14').Value = '  -2.24%  '

$ws.Range('D15').Value = '0.6746'
$ws.Range('E15').Value = '  -3.78%  '

$ws.Range('D16').Value = '297.24'
$ws.Range('E16').Value = '  +8.60%  '

$ws.Range('D17').Value = '30.538.02'
$ws.Range('E17').Value = '  -1.00%  '

$ws.Range('D18').Value = '0.9975'
$ws.Range('E18').Value = '  -0.34%  '

$ws.Range('D19').Value = '0.000007489'
$ws.Range('E19').Value = '  -2.95%  '

$ws.Range('D20').Value = '12.69'
$ws.Range('E20').Value = '  -2.81%  '

$ws.Range('D21').Value = '5.415'
$ws.Range('E21').Value = '  -2.88%  '

$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '0.9966'
$ws.Range('E22').Value = '  -0.43%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '6.317'
$ws.Range('E23').Value = '  -3.31%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '9.410'
$ws.Range('E24').Value = '  -4.48%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '167.14'
$ws.Range('E25').Value = '  +0.65%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '19.69'
$ws.Range('E26').Value = '  +0.63%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.097'
$ws.Range('E27').Value = '  -2.73%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.386'
$ws.Range('E28').Value = '  -0.38%  '

$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '0.09939'
$ws.Range('E29').Value = '  -4.41%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '4.575'
$ws.Range('E30').Value = '  +0.38%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.516'
$ws.Range('E31').Value = '  -2.82%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.262'
$ws.Range('E32').Value = '  -2.52%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.04750'
$ws.Range('E33').Value = '  -2.79%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.7283'
$ws.Range('E34').Value = '  -4.34%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.115'
$ws.Range('E35').Value = '  -3.02%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.711'
$ws.Range('E36').Value = '  -0.76%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01915'
$ws.Range('E37').Value = '  -4.68%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.618'
$ws.Range('E38').Value = '  -1.51%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '6.350'
$ws.Range('E39').Value = '  -2.36%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '74.78'
$ws.Range('E40').Value = '  -5.17%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '1.966'
$ws.Range('E41').Value = '  -5.76%  '

$ws.Range('E42').Value = '  -4.24%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '106.18'
$ws.Range('E43').Value = '  -1.48%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4273'
$ws.Range('E44').Value = '  -3.65%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9975'
$ws.Range('E45').Value = '  -0.32%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.447'
$ws.Range('E46').Value = '  -3.52%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '965.67'
$ws.Range('E47').Value = '  -3.51%  '

$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1205'
$ws.Range('E48').Value = '  -3.52%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '34.70'
$ws.Range('E49').Value = '  -4.17%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.764'
$ws.Range('E50').Value = '  -4.79%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05800'
$ws.Range('E51').Value = '  +0.29%  '
